# statistic.xlsx — update the "Insertion Sort" benchmark numbers on Лист1
# (table/chart #1, rows 2 and 4: random:comparison / sorted:comparison)
# to the freshly re-measured values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 2 = "random:comparison"
$ws.Range("C2").Value = 248341
$ws.Range("D2").Value = 991251
$ws.Range("E2").Value = 2228028
$ws.Range("F2").Value = 24946988
$ws.Range("G2").Value = 155170265
$ws.Range("H2").Value = 620681656
$ws.Range("I2").Value = 2482726624

# Row 4 = "sorted:comparison"
$ws.Range("C4").Value = 999
$ws.Range("D4").Value = 1999
$ws.Range("E4").Value = 2999
$ws.Range("F4").Value = 9999
$ws.Range("G4").Value = 24999
$ws.Range("H4").Value = 49999
$ws.Range("I4").Value = 99999

# Match the author's final active-cell selection on the sheet.
$ws.Range("I6").Select()
